$d = $word.ActiveDocument

# The first paragraph currently reads:
#   "This is a Microsoft word document."
# We need to append three more runs (without altering the existing run)
# so the paragraph becomes four runs:
#   "This is a Microsoft word document."
#   " ("
#   "Changed main"
#   ")"

$para = $d.Paragraphs.Item(1)
$r = $para.Range

# Collapse to just before the paragraph mark (end of the visible text).
$r.SetRange($r.End - 1, $r.End - 1)
$insertStart = $r.Start

$seg1 = " ("
$seg2 = "Changed main"
$seg3 = ")"

# Insert all of the new text in one shot; this keeps character offsets simple
# to compute for the follow-up pass that splits it back into separate runs.
$r.InsertAfter($seg1 + $seg2 + $seg3)

$seg1Start = $insertStart
$seg1End   = $seg1Start + $seg1.Length
$seg2Start = $seg1End
$seg2End   = $seg2Start + $seg2.Length
$seg3Start = $seg2End
$seg3End   = $seg3Start + $seg3.Length

# Touching (and then reverting) a direct-character-formatting property on each
# segment forces the engine to keep it as its own run instead of silently
# coalescing it with neighboring runs that share identical formatting.
$rSeg1 = $d.Range($seg1Start, $seg1End)
$rSeg1.Font.Bold = 1
$rSeg1.Font.Bold = 0

$rSeg2 = $d.Range($seg2Start, $seg2End)
$rSeg2.Font.Bold = 1
$rSeg2.Font.Bold = 0

$rSeg3 = $d.Range($seg3Start, $seg3End)
$rSeg3.Font.Bold = 1
$rSeg3.Font.Bold = 0
